$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.298"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'3.397"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.367"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8153"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9628"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1404"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07407"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03140"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03048"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09290"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.588"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001603"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04699"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005769"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006427"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005056"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001032"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.123"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3254"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Value = "'0.0003099"
$ws.Range("D28").Style = "Normal"
$ws.Range("D41").Value = "'0.007032"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1049"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002847"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007809"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005816"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.0005499"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.6799"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.1540"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01010"
$ws.Range("D51").Style = "Normal"
